$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 147, shifting rows 147:213 down to 148:214.
$ws.Rows("147").Insert()

# Populate the newly inserted row 147 with the new record.
$ws.Range("A147").Value = 3
$ws.Range("B147").Value = "Femacal de La Calera"
$ws.Range("C147").Value = "Coquimbo"
$ws.Range("D147").Value = "10/20/2021"
$ws.Range("E147").Value = 5
$ws.Range("F147").Value = 100112009
$ws.Range("G147").Value = "Acelga"
$ws.Range("H147").Value = "Sin especificar"
$ws.Range("I147").Value = "Primera"
$ws.Range("J147").Value = 230
$ws.Range("K147").Value = 2000
$ws.Range("L147").Value = 2200
$ws.Range("M147").Value = 2104
$ws.Range("N147").Value = "$/docena de atados (6 kilos)"
$ws.Range("O147").Value = "Provincia de Quillota"
$ws.Range("P147").Value = 351
$ws.Range("Q147").Value = 6
$ws.Range("R147").Value = "Hortaliza"
